# Roll the GSC export window forward by one day:
#   - drop the oldest date row (2025-11-01) from the "Chart" sheet
#   - shift every remaining row up by one
#   - append a new row for the new date (2026-01-30) with its Pages count

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Remove the oldest data row (row 2 = 2025-11-01) and shift everything
# below it up by one row. This naturally shifts column C's values the
# same way the diff shows, and drops "2025-11-01" from the shared
# string table once it is no longer referenced by any cell.
$ws.Range("A2:C2").Delete(-4162)

# Append the new trailing row (row 91) for 2026-01-30.
# Force the date cell to be written as literal text (not auto-converted
# to a date serial) by formatting it as Text first, then strip the
# formatting back off so the cell ends up a plain shared string again.
$newRow = 91
$dateCell = $ws.Cells.Item($newRow, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2026-01-30"
$dateCell.ClearFormats()

$ws.Cells.Item($newRow, 2).Value = 0.0
$ws.Cells.Item($newRow, 3).Value = 28.0
